$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88, pushing existing rows 88-93 down to 89-94.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new weekly record.
$ws.Range("A88").Value = 9
$ws.Range("B88").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C88").Value = "Metropolitana"
$ws.Range("D88").Value = 45106
$ws.Range("E88").Value = 13
$ws.Range("F88").Value = 100112035
$ws.Range("G88").Value = "Bruselas (repollito)"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 52
$ws.Range("K88").Value = 17000
$ws.Range("L88").Value = 19000
$ws.Range("M88").Value = 18000
$ws.Range("N88").Value = '$/malla 15 kilos'
$ws.Range("O88").Value = "Provincia de Quillota"
$ws.Range("P88").Value = 1200
$ws.Range("Q88").Value = 15
$ws.Range("R88").Value = "Hortaliza"
